$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText($sheet, $ref, $value) {
    $range = $sheet.Range($ref)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-CellText $ws "D2" "64.462.09"
Set-CellText $ws "E2" "  -0.31%  "
Set-CellText $ws "D3" "3.135.67"
Set-CellText $ws "E3" "  -1.44%  "
Set-CellText $ws "D4" "1.00"
Set-CellText $ws "E4" "  -0.01%  "
Set-CellText $ws "D5" "572.34"
Set-CellText $ws "E5" "  -0.26%  "
Set-CellText $ws "D6" "163.27"
Set-CellText $ws "E6" "  -4.34%  "
Set-CellText $ws "D8" "0.570"
Set-CellText $ws "E8" "  -6.59%  "
Set-CellText $ws "D9" "3.153.84"
Set-CellText $ws "E9" "  -1.21%  "
Set-CellText $ws "E10" "  -2.74%  "
Set-CellText $ws "D11" "6.63"
Set-CellText $ws "E11" "  -3.17%  "
Set-CellText $ws "D12" "0.380"
Set-CellText $ws "E12" "  -3.52%  "
Set-CellText $ws "D13" "3.687.26"
Set-CellText $ws "E13" "  -1.32%  "
Set-CellText $ws "E14" "  -2.26%  "
Set-CellText $ws "D15" "64.515.76"
Set-CellText $ws "E15" "  -0.29%  "
Set-CellText $ws "D16" "24.88"
Set-CellText $ws "E16" "  -2.75%  "
Set-CellText $ws "D17" "3.141.95"
Set-CellText $ws "E17" "  -2.10%  "
Set-CellText $ws "E18" "  -2.69%  "
Set-CellText $ws "D19" "412.89"
Set-CellText $ws "E19" "  -1.58%  "
Set-CellText $ws "D20" "5.23"
Set-CellText $ws "E20" "  -2.31%  "
Set-CellText $ws "D21" "12.48"
Set-CellText $ws "E21" "  -3.97%  "
Set-CellText $ws "E22" "  -2.06%  "
Set-CellText $ws "E23" "  +0.02%  "
Set-CellText $ws "D24" "68.60"
Set-CellText $ws "E24" "  -2.58%  "
Set-CellText $ws "E25" "  -3.63%  "
Set-CellText $ws "D26" "0.193"
Set-CellText $ws "E26" "  -5.71%  "
Set-CellText $ws "E27" "  -3.72%  "
Set-CellText $ws "D28" "8.91"
Set-CellText $ws "E28" "  +0.18%  "
Set-CellText $ws "D29" "0.995"
Set-CellText $ws "E29" "  -0.29%  "
Set-CellText $ws "E30" "  +0.09%  "
Set-CellText $ws "E31" "  -2.03%  "
Set-CellText $ws "D32" "21.22"
Set-CellText $ws "E32" "  -2.94%  "
Set-CellText $ws "D33" "162.70"
Set-CellText $ws "E33" "  +4.05%  "
Set-CellText $ws "D34" "4.84"
Set-CellText $ws "E34" "  -5.15%  "
Set-CellText $ws "D35" "6.25"
Set-CellText $ws "E35" "  -2.52%  "
Set-CellText $ws "E36" "  -1.47%  "
Set-CellText $ws "E37" "  -1.67%  "
Set-CellText $ws "D38" "1.67"
Set-CellText $ws "E38" "  -2.91%  "
Set-CellText $ws "D39" "2.630.99"
Set-CellText $ws "E39" "  -3.09%  "
Set-CellText $ws "B40" "EnergySwap"
Set-CellText $ws "C40" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText $ws "D40" "23.67"
Set-CellText $ws "E40" "  -3.11%  "
Set-CellText $ws "B41" "Filecoin"
Set-CellText $ws "C41" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-CellText $ws "D41" "4.11"
Set-CellText $ws "E41" "  -3.85%  "
Set-CellText $ws "D42" "38.27"
Set-CellText $ws "E42" "  -2.34%  "
Set-CellText $ws "D43" "0.691"
Set-CellText $ws "E43" "  -4.00%  "
Set-CellText $ws "D44" "0.0613"
Set-CellText $ws "E44" "  -1.91%  "
Set-CellText $ws "D45" "5.30"
Set-CellText $ws "E45" "  -5.67%  "
Set-CellText $ws "D46" "290.08"
Set-CellText $ws "E46" "  -1.72%  "
Set-CellText $ws "E47" "  -2.14%  "
Set-CellText $ws "E48" "  -3.93%  "
Set-CellText $ws "D49" "0.997"
Set-CellText $ws "E49" "  -0.11%  "
Set-CellText $ws "D50" "0.0974"
Set-CellText $ws "E50" "  -1.82%  "
Set-CellText $ws "D51" "10.47"
